# Generate Report for Handback
# The handback for f2bfd839-ba30-4254-a8da-68a9fd3cc98c.md has completed, so
# every sheet's "Ready for handoff" / stale-handback-error row for that file
# now reflects a successful handback: status flips to "Handed back: in sync
# with en-US", the handback timestamps advance, and the stale error detail
# message is cleared.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row for f2bfd839-ba30-4254-a8da-68a9fd3cc98c.md ---
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("E3").Value = "Handed back: in sync with en-US"
$ws.Range("F3").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet: row for f2bfd839-ba30-4254-a8da-68a9fd3cc98c.md ---
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("C3").Value = "Handed back: in sync with en-US"
$ws.Range("K3").Value = "2016-08-19 12:47:32"
$ws.Range("P3").Value = ""
# Error Detail column no longer holds any long error text -> shrink it back down.
$ws.Columns.Item(16).ColumnWidth = 13.7470528738839

# --- de-de sheet: row for f2bfd839-ba30-4254-a8da-68a9fd3cc98c.md ---
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("C3").Value = "Handed back: in sync with en-US"
$ws.Range("K3").Value = "2016-08-19 12:47:39"
$ws.Range("P3").Value = ""
$ws.Columns.Item(16).ColumnWidth = 13.7470528738839
